$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.10"
$ws.Range("E2").Value = "'1.34%"
$ws.Range("D3").Value = "'36.69"
$ws.Range("E3").Value = "'3.04%"
$ws.Range("D4").Value = "'5.045"
$ws.Range("E4").Value = "'-0.28%"
$ws.Range("D5").Value = "'0.07876"
$ws.Range("E5").Value = "'0.68%"
$ws.Range("D6").Value = "'2.186"
$ws.Range("E6").Value = "'-3.81%"
$ws.Range("E7").Value = "'-0.90%"
$ws.Range("D8").Value = "'4.060"
$ws.Range("E8").Value = "'1.56%"
$ws.Range("D9").Value = "'0.9265"
$ws.Range("E9").Value = "'-0.29%"
$ws.Range("D10").Value = "'0.09942"
$ws.Range("E10").Value = "'1.55%"
$ws.Range("D11").Value = "'0.1880"
$ws.Range("E11").Value = "'3.36%"
$ws.Range("D12").Value = "'0.08678"
$ws.Range("E12").Value = "'-0.56%"
$ws.Range("D13").Value = "'0.03608"
$ws.Range("E13").Value = "'5.66%"
$ws.Range("D14").Value = "'0.09954"
$ws.Range("E14").Value = "'0.28%"
$ws.Range("D15").Value = "'0.001476"
$ws.Range("E15").Value = "'-1.03%"
$ws.Range("D16").Value = "'0.005693"
$ws.Range("E16").Value = "'0.31%"
$ws.Range("D17").Value = "'3.461"
$ws.Range("E17").Value = "'-0.69%"
$ws.Range("D18").Value = "'2.472"
$ws.Range("E18").Value = "'16.21%"
$ws.Range("D19").Value = "'0.3435"
$ws.Range("E19").Value = "'0.10%"
$ws.Range("D20").Value = "'0.1328"
$ws.Range("E20").Value = "'0.57%"
$ws.Range("D21").Value = "'4.927"
$ws.Range("E21").Value = "'8.26%"
$ws.Range("D22").Value = "'0.2203"
$ws.Range("E22").Value = "'-1.47%"
$ws.Range("D23").Value = "'0.04619"
$ws.Range("E23").Value = "'-1.15%"
$ws.Range("D24").Value = "'0.005211"
$ws.Range("E24").Value = "'16.02%"
$ws.Range("E25").Value = "'0.88%"
$ws.Range("D26").Value = "'0.0001401"
$ws.Range("E26").Value = "'7.69%"
$ws.Range("D27").Value = "'0.0002720"
$ws.Range("E27").Value = "'0.76%"
$ws.Range("D39").Value = "'0.01831"
$ws.Range("E39").Value = "'4.17%"
$ws.Range("D40").Value = "'0.04764"
$ws.Range("E40").Value = "'1.35%"
$ws.Range("D41").Value = "'0.007929"
$ws.Range("E41").Value = "'-0.11%"
$ws.Range("D42").Value = "'0.1412"
$ws.Range("E42").Value = "'-0.33%"
$ws.Range("D43").Value = "'0.007601"
$ws.Range("E43").Value = "'-11.03%"
$ws.Range("D44").Value = "'0.002133"
$ws.Range("E44").Value = "'-7.32%"
$ws.Range("D45").Value = "'0.01011"
$ws.Range("E45").Value = "'10.68%"
$ws.Range("D46").Value = "'0.00006278"
$ws.Range("E46").Value = "'2.27%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.01%"
$ws.Range("D48").Value = "'0.0005807"
$ws.Range("E48").Value = "'0.12%"
$ws.Range("D49").Value = "'36.22"
$ws.Range("E49").Value = "'814.27%"
$ws.Range("D50").Value = "'0.002692"
$ws.Range("E50").Value = "'0.05%"
$ws.Range("D51").Value = "'0.00002101"
$ws.Range("E51").Value = "'0.01%"
